$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227; existing rows 227-262 shift down to 228-263
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new record
$ws.Cells.Item(227, 1).Value = 3
$ws.Cells.Item(227, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(227, 3).Value = "Coquimbo"
$ws.Cells.Item(227, 4).NumberFormat = $ws.Cells.Item(228, 4).NumberFormat
$ws.Cells.Item(227, 4).Value = 44522
$ws.Cells.Item(227, 5).Value = 5
$ws.Cells.Item(227, 6).Value = 100112013
$ws.Cells.Item(227, 7).Value = "Alcachofa"
$ws.Cells.Item(227, 8).Value = "Española"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 19900
$ws.Cells.Item(227, 11).Value = 180
$ws.Cells.Item(227, 12).Value = 280
$ws.Cells.Item(227, 13).Value = 242
$ws.Cells.Item(227, 14).Value = "`$/unidad"
$ws.Cells.Item(227, 15).Value = "Llay Llay"
$ws.Cells.Item(227, 16).Value = 242
$ws.Cells.Item(227, 17).Value = 1
$ws.Cells.Item(227, 18).Value = "Hortaliza"
